# [FIX] dunia perseederan product customer lead
#
# This script edits /tmp/work/before.xlsx (Sheet1, a "products" seed table) to:
#  - insert 3 new rows (motorcycle-ish "Golongan 1" products) after the existing
#    3 rows of that group, shifting the remaining rows down
#  - add a new column J "url" with image urls for every data row
#  - append a brand-new final row (FWD Soul Insurance)
#  - give the new final row's price/production_cost cells a small custom font

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 3 blank rows right after row 4 (so old rows 5-10 become rows 8-13)
# ---------------------------------------------------------------------------
$ws.Range("A5:A7").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) New column J ("url") - header + existing rows 2-4
#    (Order below matches the order the source workbook's string table was
#    built in: J1, J2, J4, J3, ...)
# ---------------------------------------------------------------------------
$ws.Range("J1").Value = "url"
$ws.Range("J2").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/14/64ec4cfcac371_2023-tesla-model-x-101-1671475309.jpeg"
$ws.Range("J4").Value = "https://img.freepik.com/premium-photo/crawler-buldozer-illustration-transportation-illustration-generative-ai_710947-95.jpg"
$ws.Range("J3").Value = "https://chakrajawara.co.id/media/nwddpddh/mengenal-mesin-diesel-common-rail-tdi-dan-diesel-konvensional.png"

# ---------------------------------------------------------------------------
# 3) Fill the 3 newly inserted rows (5,6,7) - "Golongan 1" additions
#    (columns A-D, then E (new names), then F (re-used names), then G-I, then J,
#    matching the original string-table build order)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1

$ws.Range("E5").Value = "Nissan GT-R"
$ws.Range("E6").Value = "Mazda RX-7 VeilSide"
$ws.Range("E7").Value = "GSX 1000rr "

$ws.Range("F5").Value = "Mesin Diesel 100HP"
$ws.Range("F6").Value = "Motor Listrik M1"
$ws.Range("F7").Value = "Bulldozer Metal Wheel"

$ws.Range("G5").Value = 10004
$ws.Range("H5").Value = 2000000
$ws.Range("I5").Value = 2000000

$ws.Range("G6").Value = 10005
$ws.Range("H6").Value = 3000000
$ws.Range("I6").Value = 3000000

$ws.Range("G7").Value = 10006
$ws.Range("H7").Value = 1800000
$ws.Range("I7").Value = 1800000

$ws.Range("J5").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/15/64ec4d60f0f7e_2021-nissan-gt-r-2457-3-1664901335.jpeg"
$ws.Range("J6").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/13/64ec4cb1f20df_IMG-20200506-WA0004-e1588733192512.jpeg"
$ws.Range("J7").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/30/64ed5c2e310a8_maxresdefault.jpeg"

# ---------------------------------------------------------------------------
# 4) Column J for the rows that were shifted down (old rows 5-10 -> 8-13)
# ---------------------------------------------------------------------------
$ws.Range("J8").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/8/64ec4a9200e31_5fdebc94f4196.jpeg"
$ws.Range("J9").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/12/64ec4bb7ee80c_32127_24959.jpeg"
$ws.Range("J10").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/7/64ec4a5a50cf3_lampukristal.jpeg"
$ws.Range("J11").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/11/64ec4b710dd2c_Screen-Shot-2023-08-28-at-14.23.16.png"
$ws.Range("J12").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/10/64ec4b4b857b9_nationwide-mutual-insurance-company4591.jpeg"
$ws.Range("J13").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/9/64ec4b031884a_insurance.jpeg"

# ---------------------------------------------------------------------------
# 5) Brand new row 14 - "FWD Soul Insurance"
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "FWD Soul Insurance"
$ws.Range("F14").Value = "FWD Soul Insurance"
$ws.Range("G14").Value = 30004
$ws.Range("H14").Value = 17000000
$ws.Range("I14").Value = 17000000
$ws.Range("J14").Value = "https://melandas-production.s3.ap-southeast-1.amazonaws.com/21/64ec59cd9475f_HQT8RYW6SJSGMP2YJVJV-59523127.jpeg"

# give H14/I14 a small custom font (Segoe UI 10, dark grey) like the source file
$priceStyle = $wb.Styles.Add("MelandasPriceStyle")
$priceStyle.Font.Size = 10
$priceStyle.Font.Name = "Segoe UI"
$priceStyle.Font.Color = 2696481
$ws.Range("H14").Style = "MelandasPriceStyle"
$ws.Range("I14").Style = "MelandasPriceStyle"
$wb.Styles.Item("MelandasPriceStyle").Delete()
